# Update odds values on Sheet1 as per the FlashScore data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("Q2").Value = 2.03
$ws.Range("R2").Value = 1.68
$ws.Range("AC2").Value = 9

# Row 3
$ws.Range("Q3").Value = 2.3
$ws.Range("R3").Value = 1.53

# Row 4
$ws.Range("G4").Value = 1.68

# Row 5
$ws.Range("G5").Value = 2.75
$ws.Range("I5").Value = 2.4

# Row 6
$ws.Range("G6").Value = 1.38
$ws.Range("H6").Value = 4.35
$ws.Range("I6").Value = 6.5
$ws.Range("J6").Value = 1.82
$ws.Range("K6").Value = 2.42
$ws.Range("L6").Value = 6
$ws.Range("P6").Value = 4.15
$ws.Range("Q6").Value = 1.55
$ws.Range("R6").Value = 2.15
$ws.Range("W6").Value = 6.9
$ws.Range("X6").Value = 6.2
$ws.Range("AD6").Value = 7.7
$ws.Range("AH6").Value = 16
$ws.Range("AI6").Value = 35
$ws.Range("AJ6").Value = 17
$ws.Range("AK6").Value = 110
$ws.Range("AN6").Value = 3.3
$ws.Range("AQ6").Value = 16.5
$ws.Range("AT6").Value = 3.2
$ws.Range("AU6").Value = 7.9
$ws.Range("AW6").Value = 8
$ws.Range("AX6").Value = 37
$ws.Range("AY6").Value = 37

# Row 8
$ws.Range("J8").Value = 2.63
$ws.Range("Q8").Value = 1.79
$ws.Range("R8").Value = 1.94

# Row 9
$ws.Range("J9").Value = 1.95

# Row 10
$ws.Range("R10").Value = 1.63

# Row 11
$ws.Range("R11").Value = 1.54

# Row 13
$ws.Range("G13").Value = 1.39
